$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.295.81'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.874.40'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7114'
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.58'
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3111'
$ws.Range("E8").Value = '  +1.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07741'
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("E10").Value = '  +0.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08464'
$ws.Range("E11").Value = '  +2.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.871.96'
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.205'
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7114'
$ws.Range("E14").Value = '  -0.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.28'
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.293.79'
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008315'
$ws.Range("E17").Value = '  +6.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.989'
$ws.Range("E18").Value = '  +2.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.72'
$ws.Range("E19").Value = '  -0.26%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.21'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.122.39'
$ws.Range("E21").Value = '  +0.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.818'
$ws.Range("E23").Value = '  -1.85%  '
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1609'
$ws.Range("E25").Value = '  +2.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.13'
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.021'
$ws.Range("E27").Value = '  +1.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.50'
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.512'
$ws.Range("E29").Value = '  +1.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.408'
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.319'
$ws.Range("E31").Value = '  +5.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.261'
$ws.Range("E32").Value = '  -4.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05259'
$ws.Range("E33").Value = '  +1.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.922'
$ws.Range("E34").Value = '  +1.11%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7445'
$ws.Range("E36").Value = '  +2.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.680'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  +0.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.713'
$ws.Range("E39").Value = '  +1.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.169.72'
$ws.Range("E40").Value = '  +2.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.356'
$ws.Range("E41").Value = '  +4.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.95'
$ws.Range("E42").Value = '  +1.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8850'
$ws.Range("E43").Value = '  -1.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '106.63'
$ws.Range("E44").Value = '  +4.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9997'
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.018.71'
$ws.Range("E46").Value = '  +0.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.809'
$ws.Range("E47").Value = '  +2.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5204'
$ws.Range("E48").Value = '  -1.19%  '
$ws.Range("E49").Value = '  +0.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.373'
$ws.Range("E50").Value = '  +0.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4301'
$ws.Range("E51").Value = '  +1.32%  '
